{"js": "// Replace the menu entry text \"RPC Explorer\" with \"Insight Explorer\".\nconst results = context.document.body.search(\"RPC Explorer\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Insight Explorer\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the menu entry text \"RPC Explorer\" with \"Insight Explorer\".\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.Text = \"RPC Explorer\"\n$find.Replacement.Text = \"Insight Explorer\"\n$find.Execute([ref]\"RPC Explorer\", $false, $false, $false, $false, $false, $true, 1, $false, \"Insight Explorer\", 2)\n"}
